# [UC] Update function points
#
# The "Average" row (row 72) that used to hold the average minutes/FP
# figure is removed; the regression equation now divides Function Points
# by a fixed constant (0.0508) stored in D84, and the descriptive text
# moves from E82 to D82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: drop the old "Average" label + SUM()/6 formula -------------
$ws.Range("C72:D72").Clear()

# --- Row 82/84: relabel the regression description, add the divisor ----
$ws.Range("E82").Clear()
$ws.Range("D82").Value = "(Time Spent) = (Function Points) / 0.0508"
$ws.Range("D84").Value = 0.0508

# --- New Use Cases table (rows 75:77): use the new divisor cell --------
$ws.Range("B75").Formula = '=C75/$D$84'
$ws.Range("B76").Formula = '=C76/$D$84'
$ws.Range("B77").Formula = '=C77/$D$84'

# --- Refresh the on-screen selection ------------------------------------
$ws.Range("C80").Select() | Out-Null

# --- Chart: trendline now shows its equation (intercept pinned at 0) ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$tl1 = $s1.Trendlines().Item(1)
$tl1.DisplayRSquared = $false
$tl1.DisplayEquation = $true
$tl1.Intercept = 0

# --- Chart: clear the explicit data-label position on the 2nd series ---
$dlbls = $chart.SeriesCollection().Item(2).DataLabels()
$dlbls.Position = $null
